$wb = $excel.ActiveWorkbook

# The new "Italy" market sheet mirrors the existing "Germany" sheet's
# layout/styling (same repeater list, same column widths/best-fit), so
# build it the way a user would: copy Germany to the end of the tab
# strip, rename it, then edit the market-specific cells.
$germany = $wb.Worksheets.Item("Germany")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$germany.Copy($null, $lastSheet) | Out-Null

$italy = $wb.Worksheets.Item($wb.Worksheets.Count)
$italy.Name = "Italy"

# Market-specific values.
$italy.Range("B2").Value = "Italy Market"
$italy.Range("B4").Value = "NGC-3145/T2221/T2223"

# Column B changed content, so its width needs to be re-fitted; column A
# kept the exact same repeater text as Germany so its best-fit width is
# already correct from the copy.
$italy.Columns.Item(2).AutoFit() | Out-Null

# Row 4 grew to fit the longer reference text.
$italy.Rows.Item(4).RowHeight = 28.8

# Match the recorded selection/active cell on the new sheet.
$italy.Range("B4").Select()

$italy.Activate()
